$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.11%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.05%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.05%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07631"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.58%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.943"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-15.00%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.826"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.59%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.781"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.88%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9162"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.18%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1750"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.65%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.06%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08550"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.85%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03143"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.26%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09995"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001515"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.32%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005726"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.93%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.007498"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,116.77%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.464"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.03%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3341"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.50%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1325"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.87%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.275"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1991"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "9.65%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04513"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.73%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001221"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.31%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004392"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.70%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01705"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.73%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04668"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.56%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007471"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.55%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.95%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.59%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01055"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.46%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006257"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.61%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.20%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-62.42%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8232"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.39%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.20%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
